$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-point the "Index" column formula (A3:A91): ROW(Cx)-2 -> ROW(Cx)-3 ---
# Written cell-by-cell so the engine keeps existing shared-formula groups
# (A4:A67 si=0, A68:A91 si=1) intact instead of collapsing them into a new group.
for ($r = 3; $r -le 91; $r++) {
    $ws.Cells.Item($r, 1).Formula = "=IF(C$r="""",""""," + "ROW(C$r)-3)"
}

# --- 2. New ParameterTuning results filled in for rows 45-51 (U/V/W = Feature Count columns) ---
$ws.Range("U45").Value = 0
$ws.Range("V45").Value = 0
$ws.Range("W45").Value = 864

$ws.Range("U46").Value = 0
$ws.Range("V46").Value = 0
$ws.Range("W46").Value = 432

$ws.Range("U47").Value = 0
$ws.Range("V47").Value = 0
$ws.Range("W47").Value = 1728

$ws.Range("R48").Value = "NA"
$ws.Range("S48").Value = "NA"
$ws.Range("T48").Value = "NA"
$ws.Range("U48").Value = 0
$ws.Range("V48").Value = 0
$ws.Range("W48").Value = 21600

$ws.Range("L49").Value = 32
$ws.Range("U49").Value = 0
$ws.Range("V49").Value = 0
$ws.Range("W49").Value = 96

$ws.Range("U50").Value = 0
$ws.Range("V50").Value = 0
$ws.Range("W50").Value = 384

$ws.Range("U51").Value = 0
$ws.Range("V51").Value = 0
$ws.Range("W51").Value = 384

# --- 3. Update the on-screen selection / scroll position to match the new run ---
$ws.Range("D7").Select()
